$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1011.34485
$ws.Range("I15").Value = 1011.34485
$ws.Range("K15").Value = 3034.03455
$ws.Range("M15").Value = -2865.03455
$ws.Range("H55").Value = 244.86667
$ws.Range("J55").Value = 278.83334
$ws.Range("L55").Value = 278.83334
$ws.Range("N55").Value = -706.83334
$ws.Range("H98").Value = 1800.92
$ws.Range("I98").Value = 1800.92
$ws.Range("K98").Value = 1800.92
$ws.Range("M98").Value = -302.9200000000001
$ws.Range("H112").Value = 45419.043
$ws.Range("J112").Value = 1879.625
$ws.Range("L112").Value = 5638.875
$ws.Range("N112").Value = -7854.875
$ws.Range("H122").Value = 1800.92
$ws.Range("I122").Value = 1800.92
$ws.Range("K122").Value = 5402.76
$ws.Range("M122").Value = -2952.76
$ws.Range("H138").Value = 5924.0356
$ws.Range("I138").Value = 5870.8
$ws.Range("J138").Value = 5935.609
$ws.Range("K138").Value = 17612.4
$ws.Range("L138").Value = 17806.827
$ws.Range("M138").Value = -12472.4
$ws.Range("N138").Value = -28086.827

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 4753.8184
$ws.Range("I5").Value = 327.42856
$ws.Range("K5").Value = 327.42856
$ws.Range("M5").Value = -215.42856
$ws.Range("H32").Value = 126348.67
$ws.Range("I32").Value = 126348.67
$ws.Range("K32").Value = 126348.67
$ws.Range("M32").Value = -126061.67
$ws.Range("H61").Value = 2025.3939
$ws.Range("I61").Value = 1644.6
$ws.Range("J61").Value = 5833.3335
$ws.Range("K61").Value = 1644.6
$ws.Range("L61").Value = 5833.3335
$ws.Range("M61").Value = -1432.6
$ws.Range("N61").Value = -6257.3335
$ws.Range("H74").Value = 4441.4614
$ws.Range("I74").Value = 3821.0625
$ws.Range("K74").Value = 3821.0625
$ws.Range("M74").Value = -2947.0625
$ws.Range("H77").Value = 4441.4614
$ws.Range("I77").Value = 3821.0625
$ws.Range("K77").Value = 19105.3125
$ws.Range("M77").Value = -14737.3125
$ws.Range("H110").Value = 142880350
$ws.Range("I110").Value = 200002100
$ws.Range("K110").Value = 200002100
$ws.Range("M110").Value = -200000055
$ws.Range("H122").Value = 22225858
$ws.Range("I122").Value = 33336288
$ws.Range("K122").Value = 100008864
$ws.Range("M122").Value = -100006414
$ws.Range("H132").Value = 20410458
$ws.Range("I132").Value = 22729418
$ws.Range("K132").Value = 68188254
$ws.Range("M132").Value = -68185724
$ws.Range("H136").Value = 2025.3939
$ws.Range("I136").Value = 1644.6
$ws.Range("J136").Value = 5833.3335
$ws.Range("K136").Value = 4933.799999999999
$ws.Range("L136").Value = 17500.0005
$ws.Range("M136").Value = -2383.799999999999
$ws.Range("N136").Value = -22600.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 4753.8184
$ws.Range("I4").Value = 327.42856
$ws.Range("K4").Value = 327.42856
$ws.Range("M4").Value = -212.42856
$ws.Range("H20").Value = 19255.105
$ws.Range("I20").Value = 26035.77
$ws.Range("K20").Value = 26035.77
$ws.Range("M20").Value = -25788.77

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4087.4788
$ws.Range("I31").Value = 1998.6666
$ws.Range("J31").Value = 4280.2925
$ws.Range("K31").Value = 1998.6666
$ws.Range("L31").Value = 4280.2925
$ws.Range("M31").Value = -1703.6666
$ws.Range("N31").Value = -4870.2925
$ws.Range("H34").Value = 4087.4788
$ws.Range("I34").Value = 1998.6666
$ws.Range("J34").Value = 4280.2925
$ws.Range("K34").Value = 1998.6666
$ws.Range("L34").Value = 4280.2925
$ws.Range("M34").Value = -1796.6666
$ws.Range("N34").Value = -4684.2925
$ws.Range("H58").Value = 558798.9
$ws.Range("I58").Value = 1602.4286
$ws.Range("J58").Value = 913378.44
$ws.Range("K58").Value = 1602.4286
$ws.Range("L58").Value = 913378.44
$ws.Range("M58").Value = -1399.4286
$ws.Range("N58").Value = -913784.44
$ws.Range("H94").Value = 1281.8572
$ws.Range("I94").Value = 1177.5
$ws.Range("K94").Value = 1177.5
$ws.Range("M94").Value = -726.5
$ws.Range("H108").Value = 230001
$ws.Range("J108").Value = 230001
$ws.Range("L108").Value = 230001
$ws.Range("N108").Value = -237681
$ws.Range("H122").Value = 1902.1111
$ws.Range("I122").Value = 1890.625
$ws.Range("K122").Value = 5671.875
$ws.Range("M122").Value = -3221.875
$ws.Range("H134").Value = 1899.8788
$ws.Range("I134").Value = 1407.5769
$ws.Range("K134").Value = 4222.7307
$ws.Range("M134").Value = -1687.7307
$ws.Range("H136").Value = 558798.9
$ws.Range("I136").Value = 1602.4286
$ws.Range("J136").Value = 913378.44
$ws.Range("K136").Value = 4807.2858
$ws.Range("L136").Value = 2740135.32
$ws.Range("M136").Value = -2257.2858
$ws.Range("N136").Value = -2745235.32

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 289
$ws.Range("I46").Value = 79
$ws.Range("K46").Value = 237
$ws.Range("M46").Value = -146
$ws.Range("H121").Value = 92330.09
$ws.Range("J121").Value = 144654.72
$ws.Range("L121").Value = 433964.16
$ws.Range("N121").Value = -436584.16

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1166
$ws.Range("I2").Value = 962.5455
$ws.Range("K2").Value = 962.5455
$ws.Range("M2").Value = -849.5455
$ws.Range("H102").Value = 2819.25
$ws.Range("I102").Value = 1709.6666
$ws.Range("K102").Value = 1709.6666
$ws.Range("M102").Value = -87.66660000000002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3657.9644
$ws.Range("I122").Value = 2774.611
$ws.Range("J122").Value = 5248
$ws.Range("K122").Value = 8323.832999999999
$ws.Range("L122").Value = 15744
$ws.Range("M122").Value = -5873.832999999999
$ws.Range("N122").Value = -20644
$ws.Range("H132").Value = 4390.727
$ws.Range("I132").Value = 2533.926
$ws.Range("K132").Value = 7601.778
$ws.Range("M132").Value = -5071.778
$ws.Range("H136").Value = 3828.7917
$ws.Range("I136").Value = 3543.2856
$ws.Range("J136").Value = 5827.3335
$ws.Range("K136").Value = 10629.8568
$ws.Range("L136").Value = 17482.0005
$ws.Range("M136").Value = -8079.856800000001
$ws.Range("N136").Value = -22582.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 21428
$ws.Range("I51").Value = 25713.334
$ws.Range("J51").Value = 15000
$ws.Range("K51").Value = 25713.334
$ws.Range("L51").Value = 15000
$ws.Range("M51").Value = -25203.334
$ws.Range("N51").Value = -16020
$ws.Range("H55").Value = 8525.667
$ws.Range("I55").Value = 5024
$ws.Range("J55").Value = 10276.5
$ws.Range("K55").Value = 5024
$ws.Range("L55").Value = 10276.5
$ws.Range("M55").Value = -4747
$ws.Range("N55").Value = -10830.5
$ws.Range("H107").Value = 880
$ws.Range("I107").Value = 905.6667
$ws.Range("K107").Value = 2717.0001
$ws.Range("M107").Value = -797.0001000000002
$ws.Range("H132").Value = 301813.06
$ws.Range("I132").Value = 428385.4
$ws.Range("K132").Value = 1285156.2
$ws.Range("M132").Value = -1282626.2
$ws.Range("H133").Value = 195357.5
$ws.Range("J133").Value = 195357.5
$ws.Range("L133").Value = 195357.5
$ws.Range("N133").Value = -205477.5
$ws.Range("H136").Value = 2682.861
$ws.Range("I136").Value = 1673.5186
$ws.Range("J136").Value = 5710.8887
$ws.Range("K136").Value = 5020.5558
$ws.Range("L136").Value = 17132.6661
$ws.Range("M136").Value = -2470.5558
$ws.Range("N136").Value = -22232.6661
